# Updated Comments file with changes as of 16-7
# The EHDSDevice.identifier row is removed from the mapping; the
# MedicalDevice.Product.ProductID target now maps onto
# EHDSDeviceUse.header.identifier instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the whole row containing "EHDSDevice.identifier" (row 3).
# This shifts all subsequent rows up by one.
$ws.Rows(3).Delete()

# The row that held "EHDSDeviceUse.header.identifier" is now row 16;
# give it the value that used to belong to the deleted row.
$ws.Range("B16").Value = "MedicalDevice.Product.ProductID"

# Match the cursor/selection position recorded in the saved workbook.
$ws.Range("A27").Select()
